# Auto-update draw results: append the 2025-11-17 Pick 4 draw as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 62

# Columns A and C hold values that look like a date / a plain number
# ("2025-11-17" and "251117"). Excel would normally auto-coerce those into
# a date serial / numeric value, but the source data keeps everything as
# literal text (t="str" in the sheet XML), so force text entry with a
# leading apostrophe and then restore the default "Normal" style so no
# extra number-format style sticks to the cell.
$ws.Range("A$row").Value = "'2025-11-17"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = "Pick 4"

$ws.Range("C$row").Value = "'251117"
$ws.Range("C$row").Style = "Normal"

$ws.Range("D$row").Value = "5-9-6-5"

$ws.Range("E$row").Value = "2025-11-17T21:40:38.349+04:00"
